# Insert a new weekly price-report row before the existing row 91 (Zapallo
# italiano / Vega Monumental Concepción), pushing the former rows 91-142
# down to 92-143.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 91..142 down to 92..143, opening up a blank row 91.
$ws.Rows(91).Insert()

# Populate the newly opened row 91 with the new weekly record.
$ws.Range("A91").Value = 11
$ws.Range("B91").Value = "Vega Monumental Concepción"
$ws.Range("C91").Value = "Bíobío"
$ws.Range("D91").Value = 44777
$ws.Range("E91").Value = 8
$ws.Range("F91").Value = 100112032
$ws.Range("G91").Value = "Zapallo italiano"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 220
$ws.Range("K91").Value = 17000
$ws.Range("L91").Value = 18000
$ws.Range("M91").Value = 17545
$ws.Range("N91").Value = "$/caja 50 unidades"
$ws.Range("O91").Value = "Región de Arica y Parinacota"
$ws.Range("P91").Value = 351
$ws.Range("Q91").Value = 50
$ws.Range("R91").Value = "Hortaliza"
